# Commit after add explicit wait
# Rename the first worksheet from "Sheet1" to "TeamMembers" to match the
# updated TeamMembers.xlsx test data file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "TeamMembers"
